$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the precision of the existing last row's timestamp (A37)
$ws.Range("A37").Value = 44350.83299573264

# Append the new day's data row (38)
$ws.Range("A38").Value = 44351.82712145429
$ws.Range("B38").Value = 75366
$ws.Range("C38").Value = 63521
$ws.Range("D38").Value = 3326
$ws.Range("E38").Value = 2126
$ws.Range("F38").Value = 1499
$ws.Range("G38").Value = 19871
$ws.Range("H38").Value = 1378
$ws.Range("I38").Value = 897
$ws.Range("J38").Value = 198
